$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number (45535 -> 2024-09-01)
# that must be bumped by one day (45536 -> 2024-09-02) for every data row.
for ($r = 2; $r -le 28; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45535) {
        $cell.Value = 45536
    }
}
